$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 3) ---
$ws.Range("F3").Value = "Secondary Relay(s)"
$ws.Range("H3").Value = "Comments or Conerns"

# --- Row 4: SafetyLight1Enable ---
$ws.Range("C4").Value = "GPO"

# --- Row 5: SafetyLight2Enable ---
$ws.Range("C5").Value = "GPO"
$ws.Range("E5").Value = "G"
$ws.Range("F5").Value = "None"
$ws.Range("G5").Value = "Turns on Safety Lights 2"

# --- Row 6: AccelertationEnable ---
$ws.Range("C6").Value = "GPO"
$ws.Range("D6").Value = "Switching Mosfet (Mosfet A)"
$ws.Range("E6").Value = "D"
$ws.Range("F6").Value = "A"
$ws.Range("G6").Value = "Failsafe signal to car to allow command of acceleration "
$ws.Range("H6").Value = "Confirm Mosfet state for Enable On/Off"

# --- Row 7: Acceleration ---
$ws.Range("C7").Value = "PWM"
$ws.Range("D7").Value = "RC Filter and Non-Inverting OP Amp"
$ws.Range("E7").Value = "C"
$ws.Range("F7").Value = "A"
$ws.Range("G7").Value = "Controls Acceleration of GEM through conversion of PWM to boosted Analog signal"
$ws.Range("H7").Value = "Test and Validate OP-Amp Gain"

# --- Row 8: SteeringEnable ---
$ws.Range("C8").Value = "GPO"
$ws.Range("D8").Value = "None"
$ws.Range("E8").Value = "F"
$ws.Range("F8").Value = "E"
$ws.Range("G8").Value = "Failsafe signal to car to allow command of steering motor"

# --- Row 9: SteeringDirection ---
$ws.Range("C9").Value = "GPO"
$ws.Range("D9").Value = "None"
$ws.Range("E9").Value = "H"
$ws.Range("F9").Value = "None"
$ws.Range("G9").Value = "Controls Direction of motor rotation (CCW or CW)"
$ws.Range("H9").Value = "Confirm Relay State for CW and CCW"

# --- Row 10: SteeringPower ---
$ws.Range("C10").Value = "PWM"
$ws.Range("D10").Value = "RC Filter and Non-Inverting OP Amp"
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "None"
$ws.Range("G10").Value = "Control Duty Cycle to Steering Driver"
$ws.Range("H10").Value = "Test and Validate OP-Amp Gain"

# --- Row 11: ESTOP (new row) ---
$ws.Range("B11").Value = "ESTOP"
$ws.Range("C11").Value = "GPI"
$ws.Range("D11").Value = "Voltage Divider"
$ws.Range("E11").Value = "B"
$ws.Range("F11").Value = "A"
$ws.Range("G11").Value = "Detect State of ESTOPs and Drive 'SafetyLight Enable1' to low if low signal is sent"
$ws.Range("H11").Value = "V-Divder was inverted by mistake (To be fixed)"

# --- Column widths ---
# (Runtime quantizes ColumnWidth to the nearest 1/6-character pixel grid, so the
# input is pre-compensated by the fixed 5/6 padding offset to land on the exact
# / nearest achievable stored width.)
$ws.Columns.Item(7).ColumnWidth = 69.16666666666667
$ws.Columns.Item(8).ColumnWidth = 39.330729166666664

# --- Selection ---
$ws.Range("H16").Select()

$wb.Save()
